{"js": "// Office.js (Word JavaScript API) edit script.\n// Applies the same changes described in the target unified diff:\n//   1. Collapses the three CORE COMPETENCIES detail paragraphs into one\n//      short summary line.\n//   2. Replaces the four generic bullet points under several\n//      PROFESSIONAL EXPERIENCE roles with the specific accomplishment\n//      bullets from the diff (several roles grow from 4 to 5/6 bullets).\n//   3. Appends a new \"TECHNICAL SKILLS\" section (Heading 2) with three\n//      summary paragraphs at the end of the document.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Build a lookup of exact paragraph text -> paragraph object using the\n// first (and, for this document, only) match. Re-used after every\n// mutating pass below because paragraph collections/indices shift once\n// insert/delete calls run.\nfunction indexByText(items) {\n  const map = new Map();\n  for (const p of items) {\n    const t = p.text;\n    if (!map.has(t)) map.set(t, p);\n  }\n  return map;\n}\n\nlet byText = indexByText(paragraphs.items);\n\nfunction getPara(text) {\n  const p = byText.get(text);\n  if (!p) {\n    throw new Error(\"Could not find paragraph with text: \" + text);\n  }\n  return p;\n}\n\n// Replace a contiguous run of paragraphs (identified by their current\n// text) with a new list of texts. Re-uses the first N-old paragraphs for\n// the first N-old new strings (keeping identical formatting / no pPr),\n// deletes any leftover old paragraphs if the new list is shorter, and\n// inserts additional plain paragraphs after the last reused paragraph if\n// the new list is longer.\nfunction replaceBlock(oldTexts, newTexts) {\n  const oldParas = oldTexts.map(getPara);\n\n  const shared = Math.min(oldParas.length, newTexts.length);\n  for (let i = 0; i < shared; i++) {\n    oldParas[i].insertText(newTexts[i], \"Replace\");\n  }\n\n  if (oldParas.length > newTexts.length) {\n    // Too many old paragraphs: delete the extras.\n    for (let i = newTexts.length; i < oldParas.length; i++) {\n      oldParas[i].delete();\n    }\n  } else if (newTexts.length > oldParas.length) {\n    // Too few old paragraphs: insert the remaining new ones after the\n    // last reused paragraph (inherits the plain/\"Normal\" formatting).\n    let anchor = oldParas[oldParas.length - 1];\n    for (let i = oldParas.length; i < newTexts.length; i++) {\n      anchor = anchor.insertParagraph(newTexts[i], \"After\");\n    }\n  }\n}\n\nconst BULLET = \"\\u2022\";\n\n// 1) CORE COMPETENCIES block: 3 detail paragraphs -> 1 summary paragraph.\nreplaceBlock(\n  [\n    `Research and Analytics: Survey Methodology: Design, sampling, weighting, longitudinal analysis ${BULLET} Statistical Analysis: Regression modeling, clustering, segmentation, machine learning ${BULLET} Geospatial Analysis: Spatial clustering, boundary estimation, demographic mapping ${BULLET} Data Visualization: Tableau, PowerBI, d3.js, Matplotlib, Seaborn, choropleth mapping ${BULLET} Research Management: Team leadership, methodology design, stakeholder communication`,\n    `Programming and Development: Python: Django/GeoDjango, Flask, Pandas, PySpark, SciKit-Learn, TensorFlow ${BULLET} JVM Languages: Scala (Spark), Java, Groovy ${BULLET} Web Technologies: JavaScript, React, d3.js, PHP, HTML/CSS ${BULLET} Database Languages: SQL, T-SQL, PostgreSQL/PostGIS ${BULLET} Statistical Computing: R, SPSS, SAS, Stata`,\n    `Data Infrastructure: Cloud Platforms: AWS (EC2, RDS, S3), Google Cloud Platform, Microsoft Azure ${BULLET} Big Data: Apache Spark, PySpark, Hadoop, Snowflake, dbt ${BULLET} Databases: PostgreSQL/PostGIS, MySQL, Oracle, MongoDB, Neo4j ${BULLET} Geospatial: ESRI ArcGIS, Quantum GIS, GeoServer, OSGeo, GRASS ${BULLET} DevOps: Docker, Git, CI/CD pipelines, automated testing, version control`,\n  ],\n  [`Research and Analytics ${BULLET} Programming and Development ${BULLET} Data Infrastructure`]\n);\nawait context.sync();\nparagraphs.load(\"items/text\");\nawait context.sync();\nbyText = indexByText(paragraphs.items);\n\n// 2) RESEARCH DIRECTOR - Progressive Change Campaign Committee: 4 -> 6 bullets.\nreplaceBlock(\n  [\n    \"\\u2022 Managed critical research operations for political campaigns\",\n    \"\\u2022 Conducted comprehensive polling and demographic analysis\",\n    \"\\u2022 Developed strategic recommendations based on data analysis\",\n    \"\\u2022 Led research team in support of progressive political initiatives\",\n  ],\n  [\n    \"\\u2022 Conceived, architected, and engineered FLEEM web application using Twilio API for thousands of simultaneous phone calls\",\n    \"\\u2022 Developed IVR polling system for early quantitative research supporting Senators Martin Heinrich and Elizabeth Warren\",\n    \"\\u2022 Built tabular and graphical reporting system with Python, GeoDjango, PostGIS, and Apache webserver\",\n    \"\\u2022 Designed survey deployment system facilitating thousands of simultaneous phone surveys\",\n    \"\\u2022 Significantly increased data collection efficiency through automated calling infrastructure\",\n    \"\\u2022 Managed comprehensive research operations for progressive political initiatives and candidates\",\n  ]\n);\nawait context.sync();\nparagraphs.load(\"items/text\");\nawait context.sync();\nbyText = indexByText(paragraphs.items);\n\n// 3) SOFTWARE ENGINEER - Salsa Labs, Inc.: 4 -> 5 bullets.\nreplaceBlock(\n  [\n    \"\\u2022 Developed software solutions for political campaigns and advocacy groups\",\n    \"\\u2022 Built web applications for voter engagement and campaign management\",\n    \"\\u2022 Integrated third-party APIs and data sources for campaign tools\",\n    \"\\u2022 Collaborated with political strategists to translate requirements into technical solutions\",\n  ],\n  [\n    \"\\u2022 Maintained and extended entire geospatial analysis and reporting tools for Java-based CRM system\",\n    \"\\u2022 Developed custom tile server for Web Map Service (WMS) integration using GeoTools and OpenLayers\",\n    \"\\u2022 Built geospatial analysis capabilities using Java, JavaScript, MySQL, and TileMill\",\n    \"\\u2022 Integrated mapping and visualization tools for political campaign data analysis\",\n    \"\\u2022 Collaborated with political strategists to translate geospatial requirements into technical solutions\",\n  ]\n);\nawait context.sync();\nparagraphs.load(\"items/text\");\nawait context.sync();\nbyText = indexByText(paragraphs.items);\n\n// 4) INTERIM TECHNOLOGY MANAGER - The Praxis Project: 4 -> 6 bullets.\nreplaceBlock(\n  [\n    \"\\u2022 Integrated technology solutions within organizational frameworks for social justice organizations\",\n    \"\\u2022 Developed data management systems for community organizing efforts\",\n    \"\\u2022 Provided technical training and support to nonprofit staff\",\n    \"\\u2022 Built custom applications for community engagement and advocacy\",\n  ],\n  [\n    \"\\u2022 Assisted in search for full-time CTO while performing all programmatic technology roles for multi-million dollar organization\",\n    \"\\u2022 Made all technology decisions and practices for massive multinational non-governmental organization\",\n    \"\\u2022 Wrote comprehensive frameworks for internal and external technology audits\",\n    \"\\u2022 Trained beneficiaries on spatial and Census data analysis for public health research\",\n    \"\\u2022 Trained NGO staff in web development using Drupal, PHP, and MySQL\",\n    \"\\u2022 Managed technology infrastructure supporting community health initiatives across multiple countries\",\n  ]\n);\nawait context.sync();\nparagraphs.load(\"items/text\");\nawait context.sync();\nbyText = indexByText(paragraphs.items);\n\n// 5) PROGRAMMER - Lake Research Partners: 4 -> 6 bullets.\nreplaceBlock(\n  [\n    \"\\u2022 Developed data analysis tools for political polling and research\",\n    \"\\u2022 Built statistical models for voter behavior analysis\",\n    \"\\u2022 Created data visualization tools for research presentations\",\n    \"\\u2022 Supported senior researchers with technical analysis and reporting\",\n  ],\n  [\n    \"\\u2022 Built the first collaborative and multi-actor contributed poll of polls used by the Democratic Party\",\n    \"\\u2022 Developed system that later became the Polling Consortium Database at The Analyst Institute\",\n    \"\\u2022 Worked on all aspects of questionnaire design, sampling, reporting and analysis for Congressional, Senate and Presidential elections\",\n    \"\\u2022 Conducted statistical modeling and analysis using SPSS, ArcGIS, Quantum GIS, GRASS, Stata, OSCAR, PostgreSQL, PostGIS, and Oracle\",\n    \"\\u2022 Pioneered integration of advanced mapping techniques into standard reports including choropleths and hexagonal grid maps\",\n    \"\\u2022 Developed innovative approaches to visualizing demographic and market data for enhanced client understanding\",\n  ]\n);\nawait context.sync();\nparagraphs.load(\"items/text\");\nawait context.sync();\nbyText = indexByText(paragraphs.items);\n\n// 6) FIELD DIRECTOR - The Feldman Group: 4 -> 6 bullets.\nreplaceBlock(\n  [\n    \"\\u2022 Managed field operations for political campaigns and research projects\",\n    \"\\u2022 Developed data collection and management systems for field work\",\n    \"\\u2022 Trained field staff on data collection protocols and quality control\",\n    \"\\u2022 Analyzed field data to inform campaign strategy and research findings\",\n  ],\n  [\n    \"\\u2022 Administered all quantitative and qualitative research operations ensuring reporting accuracy\",\n    \"\\u2022 Managed comprehensive survey fielding for multi-million dollar research firm\",\n    \"\\u2022 Developed and implemented data warehousing solutions for efficient storage and retrieval of research findings\",\n    \"\\u2022 Created custom reports and data visualizations based on specific client requirements\",\n    \"\\u2022 Introduced mapping and geospatial analysis into standard reporting procedures\",\n    \"\\u2022 Enhanced value of research deliverables through advanced analytical techniques using SPSS, OSCAR, PHP, and MySQL\",\n  ]\n);\nawait context.sync();\nparagraphs.load(\"items/text\");\nawait context.sync();\nbyText = indexByText(paragraphs.items);\n\n// 7) Append the new \"TECHNICAL SKILLS\" section at the very end of the body.\nconst lastPara = getPara(\n  \"\\u2022 Redistricting analysis used in court cases with rigorous methodology and expert testimony\"\n);\n\n// Insert the three plain body paragraphs first (chained off the existing\n// plain/\"Normal\" last paragraph so they inherit the same formatting with\n// no explicit paragraph style), then insert the Heading 2 title before\n// them so the final order is Heading -> para -> para -> para.\nconst skillsPara3 = lastPara.insertParagraph(\n  \"DATA INFRASTRUCTURE Cloud Platforms; Big Data; Databases; Geospatial; DevOps\",\n  \"After\"\n);\nconst skillsPara2 = skillsPara3.insertParagraph(\n  \"PROGRAMMING AND DEVELOPMENT Python; JVM Languages; Web Technologies; Database Languages; Statistical Computing\",\n  \"Before\"\n);\nconst skillsPara1 = skillsPara2.insertParagraph(\n  \"RESEARCH AND ANALYTICS Survey Methodology; Statistical Analysis; Geospatial Analysis; Data Visualization; Research Management\",\n  \"Before\"\n);\nconst heading = skillsPara1.insertParagraph(\"TECHNICAL SKILLS\", \"Before\");\nheading.style = \"Heading 2\";\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell) edit script.\n# Applies the same changes described in the target unified diff:\n#   1. Collapses the three CORE COMPETENCIES detail paragraphs into one\n#      short summary line.\n#   2. Replaces the four generic bullet points under several\n#      PROFESSIONAL EXPERIENCE roles with the specific accomplishment\n#      bullets from the diff (several roles grow from 4 to 5/6 bullets).\n#   3. Appends a new \"TECHNICAL SKILLS\" section (Heading 2) with three\n#      summary paragraphs at the end of the document.\n\nfunction Get-ParaIndexByText {\n    param($Doc, [string]$Text)\n    $count = $Doc.Paragraphs.Count\n    for ($i = 1; $i -le $count; $i++) {\n        $p = $Doc.Paragraphs.Item($i)\n        $t = $p.Range.Text.TrimEnd([char]13)\n        if ($t -eq $Text) {\n            return $i\n        }\n    }\n    throw \"Paragraph not found: $Text\"\n}\n\n# Replace a contiguous run of paragraphs (identified by their current\n# text) with a new list of texts. Reuses the first N-old paragraphs for\n# the first N-old new strings (keeping identical formatting / no pPr),\n# deletes any leftover old paragraphs if the new list is shorter, and\n# inserts additional plain paragraphs after the last reused paragraph if\n# the new list is longer.\nfunction Replace-ParagraphBlock {\n    param($Doc, [string[]]$OldTexts, [string[]]$NewTexts)\n\n    $indices = @()\n    foreach ($t in $OldTexts) {\n        $indices += (Get-ParaIndexByText $Doc $t)\n    }\n\n    $shared = [Math]::Min($indices.Count, $NewTexts.Count)\n    for ($i = 0; $i -lt $shared; $i++) {\n        $Doc.Paragraphs.Item($indices[$i]).Range.Text = $NewTexts[$i]\n    }\n\n    if ($indices.Count -gt $NewTexts.Count) {\n        # Too many old paragraphs: delete the extras (walk backwards so\n        # earlier indices stay valid).\n        for ($i = $indices.Count - 1; $i -ge $NewTexts.Count; $i--) {\n            $Doc.Paragraphs.Item($indices[$i]).Range.Delete()\n        }\n    }\n    elseif ($NewTexts.Count -gt $indices.Count) {\n        # Too few old paragraphs: insert the remaining new ones after the\n        # last reused paragraph (inherits the plain/\"Normal\" formatting).\n        $anchorIdx = $indices[$indices.Count - 1]\n        for ($i = $indices.Count; $i -lt $NewTexts.Count; $i++) {\n            $anchorPara = $Doc.Paragraphs.Item($anchorIdx)\n            $anchorPara.Range.InsertParagraphAfter()\n            $anchorIdx = $anchorIdx + 1\n            $Doc.Paragraphs.Item($anchorIdx).Range.Text = $NewTexts[$i]\n        }\n    }\n}\n\n$d = $word.ActiveDocument\n\n# 1) CORE COMPETENCIES block: 3 detail paragraphs -> 1 summary paragraph.\nReplace-ParagraphBlock $d @(\n    \"Research and Analytics: Survey Methodology: Design, sampling, weighting, longitudinal analysis \u2022 Statistical Analysis: Regression modeling, clustering, segmentation, machine learning \u2022 Geospatial Analysis: Spatial clustering, boundary estimation, demographic mapping \u2022 Data Visualization: Tableau, PowerBI, d3.js, Matplotlib, Seaborn, choropleth mapping \u2022 Research Management: Team leadership, methodology design, stakeholder communication\",\n    \"Programming and Development: Python: Django/GeoDjango, Flask, Pandas, PySpark, SciKit-Learn, TensorFlow \u2022 JVM Languages: Scala (Spark), Java, Groovy \u2022 Web Technologies: JavaScript, React, d3.js, PHP, HTML/CSS \u2022 Database Languages: SQL, T-SQL, PostgreSQL/PostGIS \u2022 Statistical Computing: R, SPSS, SAS, Stata\",\n    \"Data Infrastructure: Cloud Platforms: AWS (EC2, RDS, S3), Google Cloud Platform, Microsoft Azure \u2022 Big Data: Apache Spark, PySpark, Hadoop, Snowflake, dbt \u2022 Databases: PostgreSQL/PostGIS, MySQL, Oracle, MongoDB, Neo4j \u2022 Geospatial: ESRI ArcGIS, Quantum GIS, GeoServer, OSGeo, GRASS \u2022 DevOps: Docker, Git, CI/CD pipelines, automated testing, version control\"\n) @(\n    \"Research and Analytics \u2022 Programming and Development \u2022 Data Infrastructure\"\n)\n\n# 2) RESEARCH DIRECTOR - Progressive Change Campaign Committee: 4 -> 6 bullets.\nReplace-ParagraphBlock $d @(\n    \"\u2022 Managed critical research operations for political campaigns\",\n    \"\u2022 Conducted comprehensive polling and demographic analysis\",\n    \"\u2022 Developed strategic recommendations based on data analysis\",\n    \"\u2022 Led research team in support of progressive political initiatives\"\n) @(\n    \"\u2022 Conceived, architected, and engineered FLEEM web application using Twilio API for thousands of simultaneous phone calls\",\n    \"\u2022 Developed IVR polling system for early quantitative research supporting Senators Martin Heinrich and Elizabeth Warren\",\n    \"\u2022 Built tabular and graphical reporting system with Python, GeoDjango, PostGIS, and Apache webserver\",\n    \"\u2022 Designed survey deployment system facilitating thousands of simultaneous phone surveys\",\n    \"\u2022 Significantly increased data collection efficiency through automated calling infrastructure\",\n    \"\u2022 Managed comprehensive research operations for progressive political initiatives and candidates\"\n)\n\n# 3) SOFTWARE ENGINEER - Salsa Labs, Inc.: 4 -> 5 bullets.\nReplace-ParagraphBlock $d @(\n    \"\u2022 Developed software solutions for political campaigns and advocacy groups\",\n    \"\u2022 Built web applications for voter engagement and campaign management\",\n    \"\u2022 Integrated third-party APIs and data sources for campaign tools\",\n    \"\u2022 Collaborated with political strategists to translate requirements into technical solutions\"\n) @(\n    \"\u2022 Maintained and extended entire geospatial analysis and reporting tools for Java-based CRM system\",\n    \"\u2022 Developed custom tile server for Web Map Service (WMS) integration using GeoTools and OpenLayers\",\n    \"\u2022 Built geospatial analysis capabilities using Java, JavaScript, MySQL, and TileMill\",\n    \"\u2022 Integrated mapping and visualization tools for political campaign data analysis\",\n    \"\u2022 Collaborated with political strategists to translate geospatial requirements into technical solutions\"\n)\n\n# 4) INTERIM TECHNOLOGY MANAGER - The Praxis Project: 4 -> 6 bullets.\nReplace-ParagraphBlock $d @(\n    \"\u2022 Integrated technology solutions within organizational frameworks for social justice organizations\",\n    \"\u2022 Developed data management systems for community organizing efforts\",\n    \"\u2022 Provided technical training and support to nonprofit staff\",\n    \"\u2022 Built custom applications for community engagement and advocacy\"\n) @(\n    \"\u2022 Assisted in search for full-time CTO while performing all programmatic technology roles for multi-million dollar organization\",\n    \"\u2022 Made all technology decisions and practices for massive multinational non-governmental organization\",\n    \"\u2022 Wrote comprehensive frameworks for internal and external technology audits\",\n    \"\u2022 Trained beneficiaries on spatial and Census data analysis for public health research\",\n    \"\u2022 Trained NGO staff in web development using Drupal, PHP, and MySQL\",\n    \"\u2022 Managed technology infrastructure supporting community health initiatives across multiple countries\"\n)\n\n# 5) PROGRAMMER - Lake Research Partners: 4 -> 6 bullets.\nReplace-ParagraphBlock $d @(\n    \"\u2022 Developed data analysis tools for political polling and research\",\n    \"\u2022 Built statistical models for voter behavior analysis\",\n    \"\u2022 Created data visualization tools for research presentations\",\n    \"\u2022 Supported senior researchers with technical analysis and reporting\"\n) @(\n    \"\u2022 Built the first collaborative and multi-actor contributed poll of polls used by the Democratic Party\",\n    \"\u2022 Developed system that later became the Polling Consortium Database at The Analyst Institute\",\n    \"\u2022 Worked on all aspects of questionnaire design, sampling, reporting and analysis for Congressional, Senate and Presidential elections\",\n    \"\u2022 Conducted statistical modeling and analysis using SPSS, ArcGIS, Quantum GIS, GRASS, Stata, OSCAR, PostgreSQL, PostGIS, and Oracle\",\n    \"\u2022 Pioneered integration of advanced mapping techniques into standard reports including choropleths and hexagonal grid maps\",\n    \"\u2022 Developed innovative approaches to visualizing demographic and market data for enhanced client understanding\"\n)\n\n# 6) FIELD DIRECTOR - The Feldman Group: 4 -> 6 bullets.\nReplace-ParagraphBlock $d @(\n    \"\u2022 Managed field operations for political campaigns and research projects\",\n    \"\u2022 Developed data collection and management systems for field work\",\n    \"\u2022 Trained field staff on data collection protocols and quality control\",\n    \"\u2022 Analyzed field data to inform campaign strategy and research findings\"\n) @(\n    \"\u2022 Administered all quantitative and qualitative research operations ensuring reporting accuracy\",\n    \"\u2022 Managed comprehensive survey fielding for multi-million dollar research firm\",\n    \"\u2022 Developed and implemented data warehousing solutions for efficient storage and retrieval of research findings\",\n    \"\u2022 Created custom reports and data visualizations based on specific client requirements\",\n    \"\u2022 Introduced mapping and geospatial analysis into standard reporting procedures\",\n    \"\u2022 Enhanced value of research deliverables through advanced analytical techniques using SPSS, OSCAR, PHP, and MySQL\"\n)\n\n# 7) Append the new \"TECHNICAL SKILLS\" section at the very end of the body.\n# Insert the three plain body paragraphs first (chained off the existing\n# plain/\"Normal\" last paragraph so they inherit the same formatting with\n# no explicit paragraph style), then insert the Heading 2 title before\n# them so the final order is Heading -> para -> para -> para.\n$lastIdx = Get-ParaIndexByText $d \"\u2022 Redistricting analysis used in court cases with rigorous methodology and expert testimony\"\n\n$d.Paragraphs.Item($lastIdx).Range.InsertParagraphAfter()\n$idx1 = $lastIdx + 1\n$d.Paragraphs.Item($idx1).Range.Text = \"RESEARCH AND ANALYTICS Survey Methodology; Statistical Analysis; Geospatial Analysis; Data Visualization; Research Management\"\n\n$d.Paragraphs.Item($idx1).Range.InsertParagraphAfter()\n$idx2 = $idx1 + 1\n$d.Paragraphs.Item($idx2).Range.Text = \"PROGRAMMING AND DEVELOPMENT Python; JVM Languages; Web Technologies; Database Languages; Statistical Computing\"\n\n$d.Paragraphs.Item($idx2).Range.InsertParagraphAfter()\n$idx3 = $idx2 + 1\n$d.Paragraphs.Item($idx3).Range.Text = \"DATA INFRASTRUCTURE Cloud Platforms; Big Data; Databases; Geospatial; DevOps\"\n\n$d.Paragraphs.Item($idx1).Range.InsertParagraphBefore()\n$headingIdx = $idx1\n$heading = $d.Paragraphs.Item($headingIdx)\n$heading.Range.Text = \"TECHNICAL SKILLS\"\n$heading.Style = \"Heading 2\"\n"}
